$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.567.13'
$ws.Range('E2').Value = '  -1.61%  '
$ws.Range('D3').Value = '2.447.20'
$ws.Range('E3').Value = '  -2.55%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '563.43'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -2.13%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '162.70'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -2.45%  '
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('E8').Value = '  -1.73%  '
$ws.Range('E9').Value = '  -6.69%  '
$ws.Range('E10').Value = '  -1.91%  '
$ws.Range('E11').Value = '  -3.80%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.81'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -2.32%  '
$ws.Range('D13').Value = '2.898.29'
$ws.Range('E13').Value = '  -2.64%  '
$ws.Range('D14').Value = '68.405.62'
$ws.Range('E14').Value = '  -1.68%  '
$ws.Range('E15').Value = '  -4.12%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '23.67'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -4.94%  '
$ws.Range('D17').Value = '2.459.29'
$ws.Range('E17').Value = '  -2.08%  '
$ws.Range('E18').Value = '  -2.82%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '346.64'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -1.10%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.18'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -4.66%  '
$ws.Range('E21').Value = '  -2.70%  '
$ws.Range('B22').Value = 'Dai'
$ws.Range('C22').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.00'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -0.21%  '
$ws.Range('B23').Value = 'SuiNetwork'
$ws.Range('C23').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.88'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -3.62%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '68.28'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -2.76%  '
$ws.Range('E25').Value = '  -4.83%  '
$ws.Range('E26').Value = '  -2.05%  '
$ws.Range('E27').Value = '  +1.65%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.26'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -6.63%  '
$ws.Range('D29').Value = '0.0₃0839'
$ws.Range('E29').Value = '  -5.99%  '
$ws.Range('E30').Value = '  -6.81%  '
$ws.Range('E31').Value = '  -3.74%  '
$ws.Range('E32').Value = '  +0.01%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '430.75'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -6.83%  '
$ws.Range('E34').Value = '  -3.01%  '
$ws.Range('B35').Value = 'Monero'
$ws.Range('C35').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '156.56'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -1.41%  '
$ws.Range('B36').Value = 'POPCAT'
$ws.Range('C36').Value = 'https://coinranking.com/coin/sLBuDEsp6+popcat-popcat'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.94'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +98.44%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '18.99'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -0.37%  '
$ws.Range('E38').Value = '  +0.00%  '
$ws.Range('E39').Value = '  -5.74%  '
$ws.Range('E40').Value = '  -3.04%  '
$ws.Range('E41').Value = '  -4.04%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '4.50'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -4.08%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.54'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -4.26%  '
$ws.Range('E44').Value = '  +0.15%  '
$ws.Range('E45').Value = '  -6.13%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '135.14'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -5.15%  '
$ws.Range('E47').Value = '  -2.92%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.491'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -5.47%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0717'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -2.24%  '
$ws.Range('E50').Value = '  -3.09%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0915'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -1.48%  '
